# BUGFIX: Title page font size became tiny when the user clicked in the
# title page -> solved.
#
# 1) The "Subtitle 2" placeholder on the "Title Slide" custom layout had
#    its box height (and hence the auto-fit font size) shrink to a sliver
#    (644842 EMU ~= 50.77pt). Restore/grow the box so normAutofit no
#    longer has to crush the font: height -> 1533525 EMU (120.75pt), and
#    keep the (practically unchanged, off-by-one-EMU) top offset in sync
#    with the canonical edit (4714874 EMU).
# 2) The cached "datetimeFigureOut" field text on the slide master and
#    every slide layout gets refreshed from 11/8/2019 to 1/7/2022 (this
#    is the date PowerPoint re-stamps into the placeholder cache whenever
#    it touches/saves the template).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

# --- 1) Fix the Title Slide layout's Subtitle placeholder box height ---
$titleLayout = $master.CustomLayouts.Item(1)
$subtitle = $null
for ($i = 1; $i -le $titleLayout.Shapes.Count; $i++) {
    $shp = $titleLayout.Shapes.Item($i)
    if ($shp.Name -eq "Subtitle 2") {
        $subtitle = $shp
        break
    }
}
if ($subtitle -ne $null) {
    # 4714874 EMU expressed in points (1 pt = 12700 EMU); chosen so the
    # float32 round-trip through the COM Top/Height setters lands back on
    # the exact target EMU value.
    $subtitle.Top = 371.24993896484375
    $subtitle.Height = 120.75
}

# --- 2) Refresh the cached date field text everywhere it appears ---
function Update-DateField($shp) {
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "11/8/2019") {
                $tr.Text = "1/7/2022"
            }
        }
    }
}

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateField $master.Shapes.Item($i)
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateField $layout.Shapes.Item($i)
    }
}
